# Daily-Crawler update: prepend the two newest trading days (2021/11/01 and
# 2021/11/02) to every sheet, pushing the existing history rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 台指期換倉成本計算  (date / month / settle / OI / amount / cost)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert 2021/11/01 first, then 2021/11/02, so 11/02 ends up on top (row 2)
# and 11/01 underneath it (row 3) - newest date first.
$ws1.Rows.Item(2).Insert()
$ws1.Range("A2:XFD2").Select()
$ws1.Cells.Item(2,1).Value = "日期：2021/11/01"
$ws1.Cells.Item(2,2).Value = "'202112"
$ws1.Cells.Item(2,3).Value = 17053
$ws1.Cells.Item(2,4).Value = 8300
$ws1.Cells.Item(2,5).Value = 2660268
$ws1.Cells.Item(2,6).Value = 16828

$ws1.Rows.Item(2).Insert()
$ws1.Range("A2:XFD2").Select()
$ws1.Cells.Item(2,1).Value = "日期：2021/11/02"
$ws1.Cells.Item(2,2).Value = "'202112"
$ws1.Cells.Item(2,3).Value = 17020
$ws1.Cells.Item(2,4).Value = 8647
$ws1.Cells.Item(2,5).Value = 5905940
$ws1.Cells.Item(2,6).Value = 16835

$ws1.Activate()

# ---------------------------------------------------------------------
# Sheet 2: 散戶多空力道  (date / retail long-short strength)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(2).Insert()
$ws2.Range("A2:XFD2").Select()
$ws2.Cells.Item(2,1).Value = "日期：2021/11/01"
$ws2.Cells.Item(2,2).Value = 0.05

$ws2.Rows.Item(2).Insert()
$ws2.Range("A2:XFD2").Select()
$ws2.Cells.Item(2,1).Value = "日期：2021/11/02"
$ws2.Cells.Item(2,2).Value = 0.08

# ---------------------------------------------------------------------
# Sheet 3: 三大法人買賣金額  (date / foreign / domestic)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Rows.Item(2).Insert()
$ws3.Cells.Item(2,1).Value = "110年11月01日"
$ws3.Cells.Item(2,2).Value = -31.96
$ws3.Cells.Item(2,3).Value = 112.79

$ws3.Rows.Item(2).Insert()
$ws3.Cells.Item(2,1).Value = "110年11月02日"
$ws3.Cells.Item(2,2).Value = 19.71
$ws3.Cells.Item(2,3).Value = -21.98

# ---------------------------------------------------------------------
# Sheet 4: 大盤多空點位  (date / overnight long-short level)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Rows.Item(2).Insert()
$ws4.Range("A2:XFD2").Select()
$ws4.Cells.Item(2,1).Value = "110年11月01日"
$ws4.Cells.Item(2,2).Value = 17080.39

$ws4.Rows.Item(2).Insert()
$ws4.Range("A2:XFD2").Select()
$ws4.Cells.Item(2,1).Value = "110年11月02日"
$ws4.Cells.Item(2,2).Value = 17133.93

# Make sheet 1 the active tab again (activeTab=0), matching the saved view.
$ws1.Activate()
$ws1.Range("A2:F2").Select()

Write-Output "daily crawler update applied"
